$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")
$ws.Activate()

# New "Alarm Normal load" detail columns next to the existing Accessories table
# (mirrors the header style used by A7:D7 and the value style used by B8:D8)
$ws.Cells.Item(7, 1).Copy()
$ws.Cells.Item(7, 14).PasteSpecial(-4122)
$ws.Cells.Item(7, 14).Value = "AlarmLoadingDetail"

$ws.Cells.Item(7, 1).Copy()
$ws.Cells.Item(7, 15).PasteSpecial(-4122)
$ws.Cells.Item(7, 15).Value = "StandbyLoadingDetail"

$ws.Cells.Item(8, 2).Copy()
$ws.Cells.Item(8, 14).PasteSpecial(-4122)
$ws.Cells.Item(8, 14).Value = "Battery Alarm (A)"

$ws.Cells.Item(8, 2).Copy()
$ws.Cells.Item(8, 15).PasteSpecial(-4122)
$ws.Cells.Item(8, 15).Value = "Battery Standby (A)"

# Scroll the view toward the new columns and select the newly added range
$excel.ActiveWindow.ScrollColumn = 11
$ws.Range("N7:O8").Select()
